$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Corrected IFRS consolidated financial figures for the CJ Logistics company_list sheet.
# Target values per row (column letter => new numeric value).
$rowData = @{
    2 = @{ "D"=45601; "E"=1671; "F"=1671; "G"=865; "H"=595; "I"=573; "J"=22; "K"=45429; "L"=22622; "M"=22807; "N"=22388; "O"=418; "P"=1141; "Q"=2598; "R"=-674; "S"=-1793; "T"=864; "U"=1734; "V"=15282; "W"=3.66; "X"=1.3; "Y"=2.57; "Z"=1.3; "AA"=99.19; "AB"=2347.36; "AC"=2510; "AD"=78.09; "AE"=128494; "AF"=1.53; "AG"=0; "AH"=0; "AI"=0; "AJ"=22812344 }
    3 = @{ "D"=50558; "E"=1866; "F"=1866; "G"=803; "H"=490; "I"=459; "J"=30; "K"=45005; "L"=21288; "M"=23717; "N"=23233; "O"=484; "P"=1141; "Q"=1719; "R"=-864; "S"=-1008; "T"=873; "U"=846; "V"=14047; "W"=3.69; "X"=0.97; "Y"=2.01; "Z"=1.08; "AA"=89.76000000000001; "AB"=2397.03; "AC"=2014; "AD"=94.84; "AE"=131810; "AF"=1.45; "AG"=0; "AH"=0; "AI"=0; "AJ"=22812344 }
    4 = @{ "D"=60819; "E"=2284; "F"=2284; "G"=910; "H"=682; "I"=558; "J"=124; "K"=55212; "L"=27819; "M"=27393; "N"=23461; "O"=3932; "P"=1141; "Q"=2084; "R"=-7055; "S"=5432; "T"=2332; "U"=-248; "V"=17962; "W"=3.76; "X"=1.12; "Y"=2.39; "Z"=1.36; "AA"=101.56; "AB"=2446.27; "AC"=2446; "AD"=73.19; "AE"=133105; "AF"=1.34; "AI"=0; "AJ"=22812344 }
    5 = @{ "D"=71104; "E"=2357; "F"=2357; "G"=884; "H"=542; "I"=468; "J"=74; "K"=62860; "L"=35262; "M"=27598; "N"=23487; "O"=4111; "P"=1141; "Q"=2464; "R"=-6896; "S"=4565; "T"=4846; "U"=-2382; "V"=22410; "W"=3.31; "X"=0.76; "Y"=1.99; "Z"=0.92; "AA"=127.77; "AB"=2453.87; "AC"=2050; "AD"=68.28; "AE"=133253; "AF"=1.05; "AI"=0; "AJ"=22812344 }
    6 = @{ "D"=92197; "E"=2427; "F"=2427; "G"=850; "H"=518; "I"=379; "K"=78767; "L"=47378; "M"=31389; "N"=26622; "P"=1141; "Q"=619; "R"=-8435; "S"=7785; "T"=5197; "U"=-4578; "V"=29212; "W"=2.63; "X"=0.5600000000000001; "Y"=1.51; "Z"=0.73; "AA"=150.94; "AB"=2508.12; "AC"=1660; "AD"=100.61; "AE"=146637; "AF"=1.14; "AI"=0; "AJ"=22812344 }
    7 = @{ "D"=103429; "E"=3010; "G"=960; "H"=564; "I"=370; "K"=87806; "L"=54005; "M"=33801; "N"=28502; "P"=1141; "Q"=6175; "R"=-5393; "S"=1257; "T"=3642; "U"=373; "W"=2.91; "X"=0.55; "Y"=1.34; "Z"=0.68; "AA"=159.77; "AC"=1621; "AD"=90.08; "AE"=156993; "AF"=0.93; "AG"=0; "AH"=0 }
    8 = @{ "D"=111957; "E"=3652; "G"=1754; "H"=1204; "I"=956; "K"=90935; "L"=55964; "M"=34971; "N"=29455; "P"=1141; "Q"=3871; "R"=-3780; "S"=426; "T"=2534; "U"=1281; "W"=3.26; "X"=1.07; "Y"=3.3; "Z"=1.35; "AA"=160.03; "AC"=4190; "AD"=34.85; "AE"=162241; "AF"=0.9; "AG"=8; "AH"=0.01; "AI"=0.18 }
    9 = @{ "D"=120128; "E"=4114; "G"=2370; "H"=1601; "I"=1268; "K"=94353; "L"=57823; "M"=36530; "N"=30717; "P"=1141; "Q"=4370; "R"=-3766; "S"=328; "T"=2417; "U"=1944; "W"=3.42; "X"=1.33; "Y"=4.22; "Z"=1.73; "AA"=158.29; "AC"=5560; "AD"=26.26; "AE"=169192; "AF"=0.86; "AG"=31; "AH"=0.02; "AI"=0.55 }
}

foreach ($row in $rowData.Keys) {
    $cols = $rowData[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# These cells are dropped entirely (not just zeroed) to match the cleaned-up rows.
$cellsToClear = @("AG4", "AH4", "AG5", "AH5", "AG6", "AH6", "AI7")
foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}
